# Update the Multi-Level Bill of Materials sheet: add unit prices (N column)
# for several parts, change the quantity of part #10 (row 18), and let the
# totals in row 32 (and the summary cells C5/C6) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Multi-Level Bill of Materials")

# Quantity change: part #10 (row 18) goes from 4 to 1
$ws.Range("E18").Value = 1

# New unit prices in column N for rows 9, 11-18 (row 10 already had a price)
$ws.Range("N9").Value = 22.8
$ws.Range("N11").Value = 7.37
$ws.Range("N12").Value = 9.99
$ws.Range("N13").Value = 22
$ws.Range("N14").Value = 17.2
$ws.Range("N15").Value = 15
$ws.Range("N16").Value = 10.13
$ws.Range("N17").Value = 3.81
$ws.Range("N18").Value = 9.99

# Force a full recalculation so the O column (E*N), row 32 totals, and the
# summary sheet cells (C5/C6 on the first sheet) pick up the new values.
$excel.CalculateFullRebuild()

# Update the active view position/selection to match the saved state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("E18").Select()
